$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.193987488746643
$ws.Range("B1").Value = 2.315817594528198
$ws.Range("C1").Value = 3.449532032012939
$ws.Range("D1").Value = 3.336933135986328
$ws.Range("E1").Value = 1.142379522323608
